$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.376.46'
$ws.Range('E2').Value = '  +1.53%  '
$ws.Range('D3').Value = '3.932.80'
$ws.Range('E3').Value = '  +0.43%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').Value = '''487.85'
$ws.Range('E5').Value = '  +3.52%  '
$ws.Range('D6').Value = '''148.95'
$ws.Range('E6').Value = '  +3.00%  '
$ws.Range('E7').Value = '  +1.04%  '
$ws.Range('D9').Value = '''0.733'
$ws.Range('E9').Value = '  -0.40%  '
$ws.Range('E10').Value = '  +2.86%  '
$ws.Range('D11').Value = '''0.0000355'
$ws.Range('E11').Value = '  +4.83%  '
$ws.Range('D12').Value = '''43.06'
$ws.Range('E12').Value = '  -0.66%  '
$ws.Range('D13').Value = '''10.75'
$ws.Range('E13').Value = '  +3.40%  '
$ws.Range('D14').Value = '4.555.39'
$ws.Range('E14').Value = '  +0.64%  '
$ws.Range('D15').Value = '''14.70'
$ws.Range('E15').Value = '  -2.71%  '
$ws.Range('D16').Value = '3.941.74'
$ws.Range('E16').Value = '  +1.35%  '
$ws.Range('E17').Value = '  -0.08%  '
$ws.Range('D18').Value = '''20.01'
$ws.Range('E18').Value = '  +0.62%  '
$ws.Range('E19').Value = '  -2.66%  '
$ws.Range('D20').Value = '68.543.56'
$ws.Range('E20').Value = '  +1.37%  '
$ws.Range('D21').Value = '''443.52'
$ws.Range('E21').Value = '  +2.43%  '
$ws.Range('D22').Value = '''3.46'
$ws.Range('E22').Value = '  +3.06%  '
$ws.Range('D23').Value = '''14.88'
$ws.Range('E23').Value = '  +1.71%  '
$ws.Range('D24').Value = '''88.69'
$ws.Range('E24').Value = '  +0.49%  '
$ws.Range('D25').Value = '''11.41'
$ws.Range('E25').Value = '  +13.10%  '
$ws.Range('D26').Value = '''10.97'
$ws.Range('E26').Value = '  +14.98%  '
$ws.Range('D27').Value = '''3.65'
$ws.Range('E27').Value = '  +2.29%  '
$ws.Range('D28').Value = '''38.84'
$ws.Range('E28').Value = '  +0.57%  '
$ws.Range('D29').Value = '''5.89'
$ws.Range('E29').Value = '  +2.46%  '
$ws.Range('D30').Value = '''713.20'
$ws.Range('E30').Value = '  -1.75%  '
$ws.Range('D31').Value = '''13.62'
$ws.Range('E31').Value = '  -0.80%  '
$ws.Range('E32').Value = '  -1.13%  '
$ws.Range('D34').Value = '0.0₃0914'
$ws.Range('E34').Value = '  +14.08%  '
$ws.Range('D35').Value = '''41.96'
$ws.Range('E35').Value = '  -2.86%  '
$ws.Range('D36').Value = '''6.10'
$ws.Range('E36').Value = '  +12.94%  '
$ws.Range('D37').Value = '''60.97'
$ws.Range('E37').Value = '  +5.97%  '
$ws.Range('E38').Value = '  -3.39%  '
$ws.Range('D39').Value = '''0.398'
$ws.Range('E39').Value = '  +18.39%  '
$ws.Range('D40').Value = '''1.00'
$ws.Range('E40').Value = '  +0.08%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').Value = '''0.0482'
$ws.Range('E41').Value = '  +1.30%  '
$ws.Range('B42').Value = 'Fetch.AI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D42').Value = '''2.89'
$ws.Range('E42').Value = '  +13.67%  '
$ws.Range('D43').Value = '''3.19'
$ws.Range('E43').Value = '  +3.69%  '
$ws.Range('D44').Value = '''2.94'
$ws.Range('E44').Value = '  +5.36%  '
$ws.Range('D45').Value = '''0.142'
$ws.Range('E45').Value = '  +0.70%  '
$ws.Range('B46').Value = 'FirstDigitalUSD'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D46').Value = '''1.00'
$ws.Range('E46').Value = '  -0.06%  '
$ws.Range('B47').Value = 'LidoDAOToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D47').Value = '''3.41'
$ws.Range('E47').Value = '  +0.42%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0353'
$ws.Range('E48').Value = '  +41.50%  '
$ws.Range('E49').Value = '  -1.08%  '
$ws.Range('D50').Value = '''145.98'
$ws.Range('E50').Value = '  +0.46%  '
$ws.Range('D51').Value = '''3.15'
$ws.Range('E51').Value = '  -0.48%  '
